# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" right before the "总计" sheet, populated
#   with the per-fund holding detail for the quarter (same layout as the
#   other quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#   持有市值(亿元)/仓位排名).
# - Prepend a "2022-Q1" summary row to the "总计" sheet (date / holding
#   count / holding value), pushing the existing rows down by one.
#
# IMPORTANT: worksheet handles returned by Worksheets.Item(..) in this
# host are positional snapshots, not durable object references. Once a
# sheet gets inserted/copied (and the tab order shifts), a handle that
# was captured *before* the shift keeps resolving by its *old* position,
# silently pointing at the wrong sheet. To stay correct we re-resolve
# sheets (by name, or by a freshly-read Worksheets.Count) right before
# every use below instead of reusing a long-lived variable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Clone the previous quarter's sheet ("2021-Q4", the tab right before
#    "总计") immediately in front of "总计". This brings along the
#    sheetPr/pageMargins/header+column-A styling for free; we overwrite
#    its data afterwards with the 2022-Q1 numbers.
# ---------------------------------------------------------------------
$totalSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)
$templateSheet.Copy($wb.Worksheets.Item($totalSheetName))

# The copy lands immediately before "总计", i.e. one slot before the end.
$newSheetIndex = $wb.Worksheets.Count - 1
$wb.Worksheets.Item($newSheetIndex).Name = "2022-Q1"

$newSheet = $wb.Worksheets.Item("2022-Q1")

# Extend the formatted (bold + bordered) index column down to the two
# extra rows this quarter needs (7 funds vs. the template's 5).
$newSheet.Range("A6").Copy()
$newSheet.Range("A7:A8").PasteSpecial(-4122)

$funds = @(
  @{code="001167"; name="金鹰科技创新股票";           size="4.03"; stock="94.55"; pct="6.05"; mv="0.2438"; rank=1}
  @{code="210009"; name="金鹰核心资源混合";           size="3.86"; stock="94.96"; pct="6.20"; mv="0.2393"; rank=1}
  @{code="162102"; name="金鹰中小盘精选混合";         size="4.60"; stock="76.52"; pct="4.90"; mv="0.2254"; rank=1}
  @{code="010663"; name="长江均衡成长混合A";          size="0.26"; stock="85.90"; pct="3.20"; mv="0.0083"; rank=9}
  @{code="673081"; name="西部利得祥运灵活配置混合A";  size="0.14"; stock="84.19"; pct="4.43"; mv="0.0062"; rank=4}
  @{code="673083"; name="西部利得祥运灵活配置混合C";  size="0.06"; stock="84.19"; pct="4.43"; mv="0.0027"; rank=4}
  @{code="010664"; name="长江均衡成长混合C";          size="0.05"; stock="85.90"; pct="3.20"; mv="0.0016"; rank=9}
)

$r = 2
foreach ($fund in $funds) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    # Leading apostrophe forces text storage so fund codes keep their
    # leading zeros and the numeric-looking ratios stay text, matching
    # the other quarterly sheets (only the index + rank columns are
    # real numbers there).
    $newSheet.Cells.Item($r, 2).Value = "'" + $fund.code
    $newSheet.Cells.Item($r, 3).Value = $fund.name
    $newSheet.Cells.Item($r, 4).Value = "'" + $fund.size
    $newSheet.Cells.Item($r, 5).Value = "'" + $fund.stock
    $newSheet.Cells.Item($r, 6).Value = "'" + $fund.pct
    $newSheet.Cells.Item($r, 7).Value = "'" + $fund.mv
    $newSheet.Cells.Item($r, 8).Value = $fund.rank
    $r = $r + 1
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. "总计": push the existing summary rows down one slot and write the
#    new 2022-Q1 row in the gap that opens up at row 2.
# ---------------------------------------------------------------------
$quarters = @(
  @{label="2022-Q1"; count=7;  mv=0.73}
  @{label="2021-Q4"; count=5;  mv=0.53}
  @{label="2021-Q3"; count=10; mv=1.55}
  @{label="2021-Q2"; count=5;  mv=0.73}
  @{label="2021-Q1"; count=4;  mv=0.55}
  @{label="2020-Q4"; count=4;  mv=0.43}
)

$totalSheet = $wb.Worksheets.Item($totalSheetName)

# Row 7 is new: copy the formatting of the current last row (row 6) down
# onto it before the values are (re)written.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$r = 2
foreach ($q in $quarters) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $q.label
    $totalSheet.Cells.Item($r, 3).Value = $q.count
    $totalSheet.Cells.Item($r, 4).Value = $q.mv
    $r = $r + 1
}

$totalSheet.Range("A1").Select()
